# Add the new user/cage row that was missing, fixing the "register and
# login" bug where two cages could not share the same number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(7, 1).Value = "orelhen"
$ws.Cells.Item(7, 2).Value = "1234567@a"
$ws.Cells.Item(7, 3).Value = 123456789
